$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.77850740466121
$ws.Range("D2").Value = 9.495261758223313
$ws.Range("E2").Value = 15.9014303887987
$ws.Range("F2").Value = 37.16850316517689
$ws.Range("G2").Value = 3.681687945001749
$ws.Range("J2").Value = 11.77521476605821
$ws.Range("K2").Value = 9.54365112532364
$ws.Range("L2").Value = 8.989686634232084
$ws.Range("M2").Value = 15.70272226544129
$ws.Range("N2").Value = 22.10015236453665
$ws.Range("O2").Value = 27.95671559586544
$ws.Range("B3").Value = 16.70564819522751
$ws.Range("D3").Value = 9.504882877313561
$ws.Range("E3").Value = 15.93409712915508
$ws.Range("F3").Value = 37.24430682563635
$ws.Range("G3").Value = 3.683510355783175
$ws.Range("J3").Value = 11.79209072448908
$ws.Range("K3").Value = 9.339716452807622
$ws.Range("L3").Value = 8.971029351479286
$ws.Range("M3").Value = 15.68556075396452
$ws.Range("N3").Value = 22.16187785978799
$ws.Range("O3").Value = 28.023904069832
$ws.Range("B4").Value = 16.6636835488175
$ws.Range("D4").Value = 9.511779550686898
$ws.Range("E4").Value = 15.95538392873389
$ws.Range("F4").Value = 37.29774537207848
$ws.Range("G4").Value = 3.684689938415761
$ws.Range("J4").Value = 11.80300796999173
$ws.Range("K4").Value = 9.213717417787315
$ws.Range("L4").Value = 8.960577142446086
$ws.Range("M4").Value = 15.67700287759534
$ws.Range("N4").Value = 22.20157529633601
$ws.Range("O4").Value = 28.07004035728721
$ws.Range("B5").Value = 16.64729304058897
$ws.Range("D5").Value = 9.514839292564659
$ws.Range("E5").Value = 15.96436831845352
$ws.Range("F5").Value = 37.32125439557799
$ws.Range("G5").Value = 3.685185916847539
$ws.Range("J5").Value = 11.80759686925247
$ws.Range("K5").Value = 9.162247312911878
$ws.Range("L5").Value = 8.956572929224478
$ws.Range("M5").Value = 15.67401650798209
$ws.Range("N5").Value = 22.21820568844815
$ws.Range("O5").Value = 28.09006800238951
$ws.Range("B6").Value = 16.64461469763065
$ws.Range("D6").Value = 9.515362432098096
$ws.Range("E6").Value = 15.96587890501543
$ws.Range("F6").Value = 37.32526263741814
$ws.Range("G6").Value = 3.6852691984798
$ws.Range("J6").Value = 11.80836732157526
$ws.Range("K6").Value = 9.153695377754763
$ws.Range("L6").Value = 8.955923522337061
$ws.Range("M6").Value = 15.6735509769104
$ws.Range("N6").Value = 22.2209945730313
$ws.Range("O6").Value = 28.09346763911914
$ws.Range("B7").Value = 16.66345960688153
$ws.Range("D7").Value = 9.511819805367329
$ws.Range("E7").Value = 15.95550383981827
$ws.Range("F7").Value = 37.29805541065526
$ws.Range("G7").Value = 3.684696565388091
$ws.Range("J7").Value = 11.80306928998382
$ws.Range("K7").Value = 9.213023680860175
$ws.Range("L7").Value = 8.960522103533078
$ws.Range("M7").Value = 15.67696056950385
$ws.Range("N7").Value = 22.20179774222799
$ws.Range("O7").Value = 28.0703054914739
$ws.Range("B8").Value = 16.75281878390405
$ws.Range("D8").Value = 9.498374122775534
$ws.Range("E8").Value = 15.91243918229513
$ws.Range("F8").Value = 37.1932085043775
$ws.Range("G8").Value = 3.682303759786345
$ws.Range("J8").Value = 11.7809185703866
$ws.Range("K8").Value = 9.473539440845931
$ws.Range("L8").Value = 8.983047078060762
$ws.Range("M8").Value = 15.69639620981849
$ws.Range("N8").Value = 22.12106298783212
$ws.Range("O8").Value = 27.97886831672484
$ws.Range("B9").Value = 16.9494068608814
$ws.Range("D9").Value = 9.479832824346484
$ws.Range("E9").Value = 15.83771026229672
$ws.Range("F9").Value = 37.04236412950846
$ws.Range("G9").Value = 3.678090295894985
$ws.Range("J9").Value = 11.74186964908907
$ws.Range("K9").Value = 9.975155935335945
$ws.Range("L9").Value = 9.035050666520505
$ws.Range("M9").Value = 15.75006140474004
$ws.Range("N9").Value = 21.97694582564857
$ws.Range("O9").Value = 27.83833895449647
$ws.Range("B10").Value = 17.10596084525239
$ws.Range("D10").Value = 9.470948769925934
$ws.Range("E10").Value = 15.78868697237625
$ws.Range("F10").Value = 36.96498192270537
$ws.Range("G10").Value = 3.675283584062934
$ws.Range("J10").Value = 11.71583081199342
$ws.Range("K10").Value = 10.33418293054735
$ws.Range("L10").Value = 9.077850658819136
$ws.Range("M10").Value = 15.79875272799007
$ws.Range("N10").Value = 21.87963529514246
$ws.Range("O10").Value = 27.75878238361203
$ws.Range("B11").Value = 17.17960764471222
$ws.Range("D11").Value = 9.467928663883015
$ws.Range("E11").Value = 15.76765208337408
$ws.Range("F11").Value = 36.93704753444532
$ws.Range("G11").Value = 3.674068840905536
$ws.Range("J11").Value = 11.70455527159366
$ws.Range("K11").Value = 10.4946980144586
$ws.Range("L11").Value = 9.098277055669227
$ws.Range("M11").Value = 15.82286411742238
$ws.Range("N11").Value = 21.83720917298587
$ws.Range("O11").Value = 27.72774243902249
$ws.Range("B12").Value = 17.2078268755682
$ws.Range("D12").Value = 9.466931216544102
$ws.Range("E12").Value = 15.75986804022215
$ws.Range("F12").Value = 36.92751448961203
$ws.Range("G12").Value = 3.673617723174646
$ws.Range("J12").Value = 11.70036703577574
$ws.Range("K12").Value = 10.55501963539482
$ws.Range("L12").Value = 9.106145749196838
$ws.Range("M12").Value = 15.83227178849095
$ws.Range("N12").Value = 21.82140686612494
$ws.Range("O12").Value = 27.71672953270146
$ws.Range("B13").Value = 17.20173494879985
$ws.Range("D13").Value = 9.467139543145485
$ws.Range("E13").Value = 15.76153641499147
$ws.Range("F13").Value = 36.92952112242527
$ws.Range("G13").Value = 3.67371448522375
$ws.Range("J13").Value = 11.70126542512826
$ws.Range("K13").Value = 10.54204972471021
$ws.Range("L13").Value = 9.104445203866002
$ws.Range("M13").Value = 15.83023342871379
$ws.Range("N13").Value = 21.82479847306078
$ws.Range("O13").Value = 27.71906838438425
$ws.Range("B14").Value = 17.18192273663226
$ws.Range("D14").Value = 9.467843677078031
$ws.Range("E14").Value = 15.76700805334002
$ws.Range("F14").Value = 36.93624230041684
$ws.Range("G14").Value = 3.674031549470427
$ws.Range("J14").Value = 11.70420907018509
$ws.Range("K14").Value = 10.49967028791282
$ws.Range("L14").Value = 9.098921761763812
$ws.Range("M14").Value = 15.82363256501451
$ws.Range("N14").Value = 21.83590383256767
$ws.Range("O14").Value = 27.72682154212912
$ws.Range("B15").Value = 17.16982970241085
$ws.Range("D15").Value = 9.46829399910111
$ws.Range("E15").Value = 15.77038319612631
$ws.Range("F15").Value = 36.94049531593505
$ws.Range("G15").Value = 3.674226915540329
$ws.Range("J15").Value = 11.70602274984746
$ws.Range("K15").Value = 10.47364984003367
$ws.Range("L15").Value = 9.095555783097197
$ws.Range("M15").Value = 15.81962529881755
$ws.Range("N15").Value = 21.8427404724678
$ws.Range("O15").Value = 27.73166712190968
$ws.Range("B16").Value = 17.1011951664739
$ws.Range("D16").Value = 9.471166636282577
$ws.Range("E16").Value = 15.79008707284365
$ws.Range("F16").Value = 36.96695375037808
$ws.Range("G16").Value = 3.675364215210205
$ws.Range("J16").Value = 11.716579128093
$ws.Range("K16").Value = 10.32363132488317
$ws.Range("L16").Value = 9.076534663885132
$ws.Range("M16").Value = 15.7972160527011
$ws.Range("N16").Value = 21.88244488518007
$ws.Range("O16").Value = 27.76091461578291
$ws.Range("B17").Value = 17.05969932426426
$ws.Range("D17").Value = 9.473190023922148
$ws.Range("E17").Value = 15.80249856890904
$ws.Range("F17").Value = 36.98504661312334
$ws.Range("G17").Value = 3.676077772190752
$ws.Range("J17").Value = 11.72320077744769
$ws.Range("K17").Value = 10.23083797737526
$ws.Range("L17").Value = 9.065108081557836
$ws.Range("M17").Value = 15.7839677094173
$ws.Range("N17").Value = 21.90727291993073
$ws.Range("O17").Value = 27.78017669181833
$ws.Range("B18").Value = 17.03606184709356
$ws.Range("D18").Value = 9.474449993360947
$ws.Range("E18").Value = 15.80975653802691
$ws.Range("F18").Value = 36.99613718532505
$ws.Range("G18").Value = 3.676494033711498
$ws.Range("J18").Value = 11.72706301453667
$ws.Range("K18").Value = 10.17720461690538
$ws.Range("L18").Value = 9.058626081668018
$ws.Range("M18").Value = 15.77653253944686
$ws.Range("N18").Value = 21.92172669523532
$ws.Range("O18").Value = 27.79174051519494
$ws.Range("B19").Value = 17.02809863544267
$ws.Range("D19").Value = 9.474893134346653
$ws.Range("E19").Value = 15.81223445493189
$ws.Range("F19").Value = 37.00000973794158
$ws.Range("G19").Value = 3.676635977459942
$ws.Range("J19").Value = 11.72837992499054
$ws.Range("K19").Value = 10.1590022781608
$ws.Range("L19").Value = 9.056447007422793
$ws.Range("M19").Value = 15.77404702385078
$ws.Range("N19").Value = 21.92665030370364
$ws.Range("O19").Value = 27.79573906965793
$ws.Range("B20").Value = 17.06409297358695
$ws.Range("D20").Value = 9.472964682275911
$ws.Range("E20").Value = 15.80116501140179
$ws.Range("F20").Value = 36.98304980005709
$ws.Range("G20").Value = 3.676001208474283
$ws.Range("J20").Value = 11.72249034262416
$ws.Range("K20").Value = 10.24074344233777
$ws.Range("L20").Value = 9.066315146008741
$ws.Range("M20").Value = 15.78535891332861
$ws.Range("N20").Value = 21.90461200110822
$ws.Range("O20").Value = 27.77807603266742
$ws.Range("B21").Value = 17.18773323679596
$ws.Range("D21").Value = 9.467632893412631
$ws.Range("E21").Value = 15.76539598213057
$ws.Range("F21").Value = 36.93423976473856
$ws.Range("G21").Value = 3.673938179376948
$ws.Range("J21").Value = 11.70334223935059
$ws.Range("K21").Value = 10.51213113708858
$ws.Range("L21").Value = 9.100540534674508
$ws.Range("M21").Value = 15.82556391499388
$ws.Range("N21").Value = 21.8326347760852
$ws.Range("O21").Value = 27.72452413067904
$ws.Range("B22").Value = 17.27045840913245
$ws.Range("D22").Value = 9.465000107812029
$ws.Range("E22").Value = 15.74307597230452
$ws.Range("F22").Value = 36.90843103438819
$ws.Range("G22").Value = 3.672641605452494
$ws.Range("J22").Value = 11.69130311984241
$ws.Range("K22").Value = 10.68678185764669
$ws.Range("L22").Value = 9.123685986636938
$ws.Range("M22").Value = 15.85345360065245
$ws.Range("N22").Value = 21.78712910085887
$ws.Range("O22").Value = 27.69384558192089
$ws.Range("B23").Value = 17.226137096441
$ws.Range("D23").Value = 9.466327556333136
$ws.Range("E23").Value = 15.75489206619297
$ws.Range("F23").Value = 36.9216483308961
$ws.Range("G23").Value = 3.673328891638453
$ws.Range("J23").Value = 11.69768525157816
$ws.Range("K23").Value = 10.5938342622224
$ws.Range("L23").Value = 9.111263026947929
$ws.Range("M23").Value = 15.83842238378576
$ws.Range("N23").Value = 21.81127621481048
$ws.Range("O23").Value = 27.70982380372863
$ws.Range("B24").Value = 17.0621059210497
$ws.Range("D24").Value = 9.473066257998374
$ws.Range("E24").Value = 15.80176753134224
$ws.Range("F24").Value = 36.9839504138625
$ws.Range("G24").Value = 3.676035804150178
$ws.Range("J24").Value = 11.72281135773481
$ws.Range("K24").Value = 10.23626606558229
$ws.Range("L24").Value = 9.065769159964633
$ws.Range("M24").Value = 15.78472938438715
$ws.Range("N24").Value = 21.90581444232068
$ws.Range("O24").Value = 27.77902421501268
$ws.Range("B25").Value = 16.89402819343861
$ws.Range("D25").Value = 9.484014198569607
$ws.Range("E25").Value = 15.8568905203319
$ws.Range("F25").Value = 37.07730168888366
$ws.Range("G25").Value = 3.679179198057293
$ws.Range("J25").Value = 11.75196618504877
$ws.Range("K25").Value = 9.840858527110145
$ws.Range("L25").Value = 9.020161977183335
$ws.Range("M25").Value = 15.73390009090563
$ws.Range("N25").Value = 22.01442173764501
$ws.Range("O25").Value = 27.87219869138022
